$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells E1:G1 need the same header style (bold, centered, bordered)
# as the existing header row - copy formatting from A1 first, then set values.
$ws.Range("A1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)

# Row 1: headers (shift C/D -> D/E values and add new E/F/G names)
$ws.Range("C1").Value = "Vikas Mahajan"
$ws.Range("D1").Value = "EMO"
$ws.Range("E1").Value = "hugh"
$ws.Range("F1").Value = "raja"
$ws.Range("G1").Value = "crime master gogo"

# Row 2: Age
$ws.Range("C2").Value = 44
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 40
$ws.Range("G2").Value = 35

# Row 3: Gender
$ws.Range("C3").Value = "M"
$ws.Range("D3").Value = "F"
$ws.Range("E3").Value = "M"
$ws.Range("F3").Value = "M"
$ws.Range("G3").Value = "M"

# Row 4: Race
$ws.Range("C4").Value = 400
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 200
$ws.Range("G4").Value = 0

# Row 5: Time
$ws.Range("B5").Value = 60
$ws.Range("C5").Value = 180
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 90
$ws.Range("F5").Value = 180
$ws.Range("G5").Value = 0

# Row 6: Past height
$ws.Range("C6").Value = 1.72
$ws.Range("D6").Value = 1.72
$ws.Range("E6").Value = 1.7
$ws.Range("F6").Value = 1.57
$ws.Range("G6").Value = 1.8

# Row 7: Current height
$ws.Range("C7").Value = 1.72
$ws.Range("D7").Value = 1.77
$ws.Range("E7").Value = 1.7
$ws.Range("F7").Value = 1.57
$ws.Range("G7").Value = 1.82

# Row 8: Past BMI
$ws.Range("C8").Value = 26.36560302866415
$ws.Range("D8").Value = 15.21092482422931
$ws.Range("E8").Value = 26.98961937716263
$ws.Range("F8").Value = 31.64428577224228
$ws.Range("G8").Value = 18.20987654320988

# Row 9: Current BMI
$ws.Range("C9").Value = 25.01352082206598
$ws.Range("D9").Value = 15.95965399470139
$ws.Range("E9").Value = 28.71972318339101
$ws.Range("F9").Value = 32.04998174368129
$ws.Range("G9").Value = 19.01944209636517
